$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 74093.664
$ws.Range("I6").Value = 279.75
$ws.Range("J6").Value = 221721.5
$ws.Range("K6").Value = 839.25
$ws.Range("L6").Value = 665164.5
$ws.Range("M6").Value = -727.25
$ws.Range("N6").Value = -665388.5
$ws.Range("H107").Value = 1462.2
$ws.Range("I107").Value = 1399.5
$ws.Range("J107").Value = 1504
$ws.Range("K107").Value = 1399.5
$ws.Range("L107").Value = 1504
$ws.Range("M107").Value = 520.5
$ws.Range("N107").Value = -5344
$ws.Range("H112").Value = 55557970
$ws.Range("J112").Value = 2709.6667
$ws.Range("L112").Value = 8129.000100000001
$ws.Range("N112").Value = -10345.0001
$ws.Range("H133").Value = 25000
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -35120
$ws.Range("H137").Value = 1925863.9
$ws.Range("I137").Value = 3449995.5
$ws.Range("J137").Value = 4132.7393
$ws.Range("K137").Value = 10349986.5
$ws.Range("L137").Value = 12398.2179
$ws.Range("M137").Value = -10347436.5
$ws.Range("N137").Value = -17498.2179
$ws.Range("H138").Value = 4548584.5
$ws.Range("I138").Value = 3496.2
$ws.Range("J138").Value = 5131288
$ws.Range("K138").Value = 10488.6
$ws.Range("L138").Value = 15393864
$ws.Range("M138").Value = -5348.599999999999
$ws.Range("N138").Value = -15404144

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1371.4286
$ws.Range("I45").Value = 1325
$ws.Range("J45").Value = 1433.3334
$ws.Range("K45").Value = 1325
$ws.Range("L45").Value = 1433.3334
$ws.Range("M45").Value = -948
$ws.Range("N45").Value = -2187.3334
$ws.Range("H61").Value = 91092790
$ws.Range("I61").Value = 111224300
$ws.Range("J61").Value = 501000
$ws.Range("K61").Value = 111224300
$ws.Range("L61").Value = 501000
$ws.Range("M61").Value = -111224088
$ws.Range("N61").Value = -501424
$ws.Range("H74").Value = 6099131
$ws.Range("I74").Value = 8799750
$ws.Range("J74").Value = 62453.766
$ws.Range("K74").Value = 8799750
$ws.Range("L74").Value = 62453.766
$ws.Range("M74").Value = -8798876
$ws.Range("N74").Value = -64201.766
$ws.Range("H77").Value = 6099131
$ws.Range("I77").Value = 8799750
$ws.Range("J77").Value = 62453.766
$ws.Range("K77").Value = 43998750
$ws.Range("L77").Value = 312268.83
$ws.Range("M77").Value = -43994382
$ws.Range("N77").Value = -321004.83
$ws.Range("H110").Value = 1514.125
$ws.Range("I110").Value = 967.7273
$ws.Range("J110").Value = 2716.2
$ws.Range("K110").Value = 967.7273
$ws.Range("L110").Value = 2716.2
$ws.Range("M110").Value = 1077.2727
$ws.Range("N110").Value = -6806.2
$ws.Range("H132").Value = 10461794
$ws.Range("I132").Value = 11652616
$ws.Range("J132").Value = 220719.8
$ws.Range("K132").Value = 34957848
$ws.Range("L132").Value = 662159.3999999999
$ws.Range("M132").Value = -34955318
$ws.Range("N132").Value = -667219.3999999999
$ws.Range("H136").Value = 91092790
$ws.Range("I136").Value = 111224300
$ws.Range("J136").Value = 501000
$ws.Range("K136").Value = 333672900
$ws.Range("L136").Value = 1503000
$ws.Range("M136").Value = -333670350
$ws.Range("N136").Value = -1508100

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2146.5715
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 2171
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2171
$ws.Range("M20").Value = -1753
$ws.Range("N20").Value = -2665
$ws.Range("H134").Value = 11911239
$ws.Range("I134").Value = 7010.2
$ws.Range("J134").Value = 111113150
$ws.Range("K134").Value = 21030.6
$ws.Range("L134").Value = 333339450
$ws.Range("M134").Value = -18495.6
$ws.Range("N134").Value = -333344520
$ws.Range("H138").Value = 32898.57
$ws.Range("J138").Value = 32898.57
$ws.Range("L138").Value = 32898.57
$ws.Range("N138").Value = -43178.57

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 26992.5
$ws.Range("J64").Value = 26992.5
$ws.Range("L64").Value = 26992.5
$ws.Range("N64").Value = -27488.5
$ws.Range("H67").Value = 26992.5
$ws.Range("J67").Value = 26992.5
$ws.Range("L67").Value = 26992.5
$ws.Range("N67").Value = -28708.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3973.3333
$ws.Range("J119").Value = 4996.6
$ws.Range("L119").Value = 14989.8
$ws.Range("N119").Value = -24665.8
$ws.Range("H120").Value = 4516.6665
$ws.Range("I120").Value = 4516.6665
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 13549.9995
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -8711.999500000002
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 976.3333
$ws.Range("I122").Value = 225.4
$ws.Range("J122").Value = 1351.8
$ws.Range("K122").Value = 2028.6
$ws.Range("L122").Value = 12166.2
$ws.Range("M122").Value = 421.3999999999999
$ws.Range("N122").Value = -17066.2
$ws.Range("H123").Value = 1826
$ws.Range("I123").Value = 843.3333
$ws.Range("J123").Value = 3300
$ws.Range("K123").Value = 2529.9999
$ws.Range("L123").Value = 9900
$ws.Range("M123").Value = -79.9998999999998
$ws.Range("N123").Value = -14800
$ws.Range("H124").Value = 995.57574
$ws.Range("I124").Value = 719.5
$ws.Range("J124").Value = 1033.6552
$ws.Range("K124").Value = 2158.5
$ws.Range("L124").Value = 3100.9656
$ws.Range("M124").Value = 2751.5
$ws.Range("N124").Value = -12920.9656
$ws.Range("H125").Value = 2781.111
$ws.Range("I125").Value = 1882.5
$ws.Range("K125").Value = 5647.5
$ws.Range("M125").Value = -727.5
$ws.Range("H131").Value = 988.61536
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 1135.2
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 3405.6
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -13485.6

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 368.66666
$ws.Range("I107").Value = 229.25
$ws.Range("J107").Value = 1065.75
$ws.Range("K107").Value = 229.25
$ws.Range("L107").Value = 1065.75
$ws.Range("M107").Value = 1690.75
$ws.Range("N107").Value = -4905.75
$ws.Range("H132").Value = 64882.72
$ws.Range("I132").Value = 40148.27
$ws.Range("J132").Value = 172065.33
$ws.Range("K132").Value = 120444.81
$ws.Range("L132").Value = 516195.99
$ws.Range("M132").Value = -117914.81
$ws.Range("N132").Value = -521255.99
$ws.Range("H133").Value = 54883.332
$ws.Range("J133").Value = 54883.332
$ws.Range("L133").Value = 54883.332
$ws.Range("N133").Value = -65003.332
$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140
$ws.Range("H140").Value = 63954.547
$ws.Range("J140").Value = 63954.547
$ws.Range("L140").Value = 63954.547
$ws.Range("N140").Value = -74314.54699999999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 32004.771
$ws.Range("I132").Value = 2384
$ws.Range("J132").Value = 131974.88
$ws.Range("K132").Value = 7152
$ws.Range("L132").Value = 395924.64
$ws.Range("M132").Value = -4622
$ws.Range("N132").Value = -400984.64
$ws.Range("H136").Value = 38876
$ws.Range("I136").Value = 23997.312
$ws.Range("J136").Value = 128148.125
$ws.Range("K136").Value = 71991.936
$ws.Range("L136").Value = 384444.375
$ws.Range("M136").Value = -69441.936
$ws.Range("N136").Value = -389544.375
$ws.Range("H137").Value = 32840
$ws.Range("I137").Value = 24800
$ws.Range("J137").Value = 34850
$ws.Range("K137").Value = 24800
$ws.Range("L137").Value = 34850
$ws.Range("M137").Value = -19700
$ws.Range("N137").Value = -45050
$ws.Range("H139").Value = 36143
$ws.Range("J139").Value = 36143
$ws.Range("L139").Value = 36143
$ws.Range("N139").Value = -46423

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 36606.035
$ws.Range("I132").Value = 33217.773
$ws.Range("K132").Value = 99653.319
$ws.Range("M132").Value = -97123.319
$ws.Range("H136").Value = 50447.977
$ws.Range("I136").Value = 37307.32
$ws.Range("K136").Value = 111921.96
$ws.Range("M136").Value = -109371.96
